$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.549.56'
$ws.Range("E2").Value = '  +1.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.489.87'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '490.32'
$ws.Range("E5").Value = '  +2.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.94'
$ws.Range("E6").Value = '  +8.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("E8").Value = '  +0.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.505.61'
$ws.Range("E9").Value = '  +0.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.79'
$ws.Range("E10").Value = '  +6.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.337'
$ws.Range("E12").Value = '  +3.36%  '

$ws.Range("E13").Value = '  +1.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.923.95'
$ws.Range("E14").Value = '  -0.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.545.03'
$ws.Range("E15").Value = '  +1.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.19'
$ws.Range("E16").Value = '  +3.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  -0.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.496.01'
$ws.Range("E18").Value = '  -0.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.54'
$ws.Range("E19").Value = '  +4.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.24'
$ws.Range("E20").Value = '  +3.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.98'
$ws.Range("E21").Value = '  +0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.92'
$ws.Range("E23").Value = '  +4.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.65'
$ws.Range("E24").Value = '  +1.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.412'
$ws.Range("E25").Value = '  +2.31%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -1.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.595.72'
$ws.Range("E28").Value = '  -0.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.66'
$ws.Range("E29").Value = '  +4.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0791'
$ws.Range("E30").Value = '  +2.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.70'
$ws.Range("E32").Value = '  +0.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.38'
$ws.Range("E33").Value = '  +1.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.52'
$ws.Range("E34").Value = '  +2.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.21'
$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.17'
$ws.Range("E36").Value = '  +5.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.78'
$ws.Range("E37").Value = '  +2.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.874'
$ws.Range("E38").Value = '  +4.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.40'
$ws.Range("E39").Value = '  +6.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.95'
$ws.Range("E40").Value = '  -1.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.51'
$ws.Range("E41").Value = '  +3.90%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.613'
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0558'
$ws.Range("E43").Value = '  +1.82%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("E44").Value = '  -0.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '267.89'
$ws.Range("E45").Value = '  +7.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.87'
$ws.Range("E46").Value = '  +10.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0931'
$ws.Range("E47").Value = '  +3.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0230'
$ws.Range("E48").Value = '  +3.72%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.23'
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.82'
$ws.Range("E50").Value = '  +2.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.896.99'
$ws.Range("E51").Value = '  -3.50%  '
